# Populate the pitch-by-pitch "Catcher's View" grid (Pitch / Choice / Result
# columns F/G/H) for each at-bat, fill in a couple of Exit Velo / Launch Angle
# values that were previously blank, normalize the "Pitch Mix" ordering, and
# fix two "Undefined" results to "Strikeout" -- this is the data needed to
# drive the new strikezone visual for hitters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (Pitch, Choice, Result)
$pitches = @{
  10 = @('FB','Swing','Foul')
  11 = @('FB','Swing','Foul')
  12 = @('CB','Take','Ball')
  13 = @('FB','Swing','In Play')

  19 = @('FB','Swing','Foul')
  20 = @('FB','Swing','Foul')
  21 = @('CB','Swing','In Play')

  28 = @('CB','Take','Strike')
  29 = @('CB','Take','Ball')
  30 = @('CB','Take','Ball')
  31 = @('FB','Take','Strike')
  32 = @('CH','Take','Ball')
  33 = @('FB','Swing','Foul')
  34 = @('FB','Swing','Foul')
  35 = @('CB','Swing','Foul')

  37 = @('CB','Take','Ball')
  38 = @('CH','Take','Ball')
  39 = @('CH','Swing','Foul')
  40 = @('CH','Take','Ball')
  41 = @('CH','Swing','In Play')

  46 = @('FB','Take','Strike')
  47 = @('SL','Take','Ball')
  48 = @('SL','Take','Ball')
  49 = @('CH','Take','Ball')
  50 = @('CH','Swing','In Play')

  61 = @('CH','Take','Strike')
  62 = @('CH','Swing','Foul')
  63 = @('CH','Swing','Strike')

  70 = @('CH','Take','Ball')
  71 = @('CH','Take','Ball')
  72 = @('CH','Take','Ball')
  73 = @('CH','Take','Strike')
  74 = @('CH','Take','Strike')
  75 = @('CH','Swing','Foul')
  76 = @('CH','Take','Strike')
}

foreach ($row in $pitches.Keys) {
  $vals = $pitches[$row]
  $ws.Range("F$row").Value = $vals[0]
  $ws.Range("G$row").Value = $vals[1]
  $ws.Range("H$row").Value = $vals[2]
}

# Exit Velo (column M, header rows of each at-bat)
$ws.Range("M10").Value = "70.73 MPH"
$ws.Range("M19").Value = "82.95 MPH"
$ws.Range("M28").Value = "93.65 MPH"
$ws.Range("M37").Value = "83.46 MPH"
$ws.Range("M46").Value = "86.78 MPH"

# Launch Angle (column M, third row of each at-bat)
$ws.Range("M12").Value = "-43.16°"
$ws.Range("M21").Value = "14.74°"
$ws.Range("M30").Value = "33.31°"
$ws.Range("M39").Value = "69.2°"
$ws.Range("M48").Value = "22.04°"

# These two Launch Angle cells go back to being truly empty (no value) rather
# than holding an empty string.
$ws.Range("M63").Value = ""
$ws.Range("M72").Value = ""

# Pitch Mix lists reordered
$ws.Range("J17").Value = "CH,CB,FB"
$ws.Range("J26").Value = "CH,CB,FB"
$ws.Range("J35").Value = "CH,CB,FB"
$ws.Range("J44").Value = "CH,CB,FB,SL"
$ws.Range("J53").Value = "CH,FB,SL"
$ws.Range("J68").Value = "CH,FB,SL"
$ws.Range("J77").Value = "CH,FB,SL"

# Fix "Undefined" results
$ws.Range("M66").Value = "Strikeout"
$ws.Range("M75").Value = "Strikeout"
